$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename a handful of people (text tweaks) ---
$ws.Range("A31").Value = "Morales-Robinson, Ana Y."
$ws.Range("A9").Value  = "Ruiz-Caceres, Gaby A"
$ws.Range("A3").Value  = "Suriel, Sal"
$ws.Range("A10").Value = "Dong, Sean"

# --- Add the "Start Date 2" / "End Date 2" header columns (D, E) ---
# Copy the existing header cell so the fill style (s=3) carries over exactly,
# then overwrite with the new label.
$ws.Range("C1").Copy($ws.Range("D1"))
$ws.Range("D1").Value = "Start Date 2"
$ws.Range("C1").Copy($ws.Range("E1"))
$ws.Range("E1").Value = "End Date 2"

# --- Row 6 (Alexis, Jennifer): insert an earlier Start Date, shifting the
#     old Start/End dates one column to the right ---
$ws.Range("C6").Copy($ws.Range("D6"))
$ws.Range("B6").Copy($ws.Range("C6"))
$ws.Range("B6").Value = 40909

# --- New row 32: Agarwala, Shelly ---
$ws.Range("A31").Copy($ws.Range("A32"))
$ws.Range("A32").Value = "Agarwala, Shelly"
$ws.Range("B31").Copy($ws.Range("B32"))
$ws.Range("B32").Value = 40909

# --- Column widths for the two new columns ---
$ws.Columns("D").ColumnWidth = 10.333333333333334
$ws.Columns("E").ColumnWidth = 9.5

# --- Selection moves to G21 ---
$ws.Range("G21").Select() | Out-Null
